$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing content (A1 "Origin CMM output")
$ws.Cells.Clear()

# Values extracted from switch statement
$ws.Range("B5").Value = 3609
$ws.Range("D7").Value = 3693
$ws.Range("F10").Value = 2063
